$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: run a no-text-change Find/Replace over a unique phrase so the
# engine's run-normalization coalesces adjacent, identically-formatted runs
# within that paragraph (mirrors the "merge runs" parts of the diff).
# ---------------------------------------------------------------------------
function Normalize-Phrase([string]$phrase) {
    $rng = $d.Content
    $found = $rng.Find.Execute($phrase, $false, $false, $false, $false, $false, $true, 1, $false, $phrase, 2)
    if (-not $found) {
        Write-Host "WARNING: phrase not found for normalization: $phrase"
    }
}

# 1. Know-Center GmbH heading line: merge "Know-Center GmbH" + " " and
#    "May 2021" + " - " + "Present" + " " runs.
Normalize-Phrase("Know-Center GmbH")
Normalize-Phrase("May 2021")

# 2. Big Ass Fans heading line: merge "Big Ass Fans" + " " and the
#    "September" + " 20" + "20" + " - " + "May 2021" + " " runs.
Normalize-Phrase("Big Ass Fans")
Normalize-Phrase("September 2020")

# 3. "Marketing and Manufacturing Engineering Internship Rotations" italic runs.
Normalize-Phrase("Marketing and Manufacturing Engineering Internship Rotations")

# 4. "Cypress Semiconductor Corporation" - merge the stray "C" run.
Normalize-Phrase("Cypress Semiconductor Corporation")

# 5. "Designed a partially automated..." - merge the split runs AND
#    change the ". " bold run split so the period becomes its own run
#    with bCs (bold-complex-script) formatting instead of being part of
#    the bold run. This needs exact run-level XML control, so replace the
#    whole paragraph's content range via InsertXML.
$full = $d.Content.Text
$pStart = $full.IndexOf("Designed a partially a")
$pEndMarker = "for this task. "
$pEndIdx = $full.IndexOf($pEndMarker, $pStart) + $pEndMarker.Length
$pRange = $d.Range($pStart, $pEndIdx)
$designedXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Designed a partially automated process in Ruby and SKILL for performing large-scale quality assurance tasks</w:t></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">These improvements saved weeks of working time for this task. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pRange.InsertXML($designedXml)

# 6. "Developed a strong foundational knowledge..." bullet - remove the
#    stray "o" character (Segoe UI Symbol run) that sits between
#    "techniques. " and "Applied", leaving a single regular space.
$full = $d.Content.Text
$idx = $full.IndexOf("techniques. o Applied")
if ($idx -ge 0) {
    # "techniques." is 11 chars; delete the space that ends the first run
    # together with the stray "o" (2 chars), leaving the trailing space of
    # the (formerly "o ") run intact.
    $delRange = $d.Range($idx + 11, $idx + 13)
    $delRange.Text = ""
} else {
    Write-Host "WARNING: stray 'o' phrase not found"
}

# 7. "Crafted visually appealing..." - merge the two runs.
Normalize-Phrase("Crafted visually appealing and descriptive documentation")

# 8. "Leveraged VLSI industry-standard..." - merge the split runs.
Normalize-Phrase("Leveraged VLSI industry-standard CAD tools")

# 9. "Directed small groups..." - merge the split runs.
Normalize-Phrase("Directed small groups to rigorously test designs")

# 10. "Proposes a novel logic family..." - merge the split runs.
Normalize-Phrase("Proposes a novel logic family")

# 11. "Recipient of the University of Kentucky..." - merge the split runs.
Normalize-Phrase("Recipient of the University of Kentucky Electrical and Computer Engineering Undergraduate Research Fellowship")
